# Update the "dSF" (column F) values for a set of rows to reflect
# a repull/recalculation of the underlying data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column F
$updates = @{
    3  = 1
    9  = 2
    11 = 3
    26 = 0
    30 = 3
    40 = 0
    46 = 0
    47 = 0
    58 = 2
    64 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}

$wb.Save()
